$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write a literal text value to a cell without letting Excel
# auto-convert number-like strings (e.g. "1.001") into numeric values,
# and without leaving a residual Text number-format style behind.
function Set-CellText($ws, $ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# Refresh cryptos list: updated Price / Volume(1h) figures, and a content
# swap of rows 34-35 (ImmutableX now ranks above HuobiToken).

Set-CellText $ws 'D2' '26.831.18'
Set-CellText $ws 'E2' '  -0.88%  '
Set-CellText $ws 'D3' '1.868.18'
Set-CellText $ws 'E3' '  +0.05%  '
Set-CellText $ws 'D4' '1.001'
Set-CellText $ws 'E4' '  +0.09%  '
Set-CellText $ws 'D5' '304.98'
Set-CellText $ws 'E5' '  -0.34%  '
Set-CellText $ws 'D6' '1.001'
Set-CellText $ws 'E6' '  +0.18%  '
Set-CellText $ws 'D7' '0.5089'
Set-CellText $ws 'E7' '  -0.83%  '
Set-CellText $ws 'D8' '0.3656'
Set-CellText $ws 'E8' '  -2.70%  '
Set-CellText $ws 'D9' '0.07182'
Set-CellText $ws 'E9' '  +0.30%  '
Set-CellText $ws 'D10' '0.8921'
Set-CellText $ws 'E10' '  +0.27%  '
Set-CellText $ws 'D11' '20.69'
Set-CellText $ws 'E11' '  +0.00%  '
Set-CellText $ws 'D12' '0.07520'
Set-CellText $ws 'E12' '  -1.08%  '
Set-CellText $ws 'D13' '1.875.72'
Set-CellText $ws 'E13' '  +0.57%  '
Set-CellText $ws 'D14' '94.84'
Set-CellText $ws 'E14' '  +5.77%  '
Set-CellText $ws 'D15' '5.220'
Set-CellText $ws 'E15' '  -1.47%  '
Set-CellText $ws 'D16' '1.001'
Set-CellText $ws 'E16' '  +0.04%  '
Set-CellText $ws 'D17' '0.000008494'
Set-CellText $ws 'E17' '  +0.25%  '
Set-CellText $ws 'D18' '14.18'
Set-CellText $ws 'E18' '  +0.77%  '
Set-CellText $ws 'D19' '1.002'
Set-CellText $ws 'E19' '  +0.20%  '
Set-CellText $ws 'D20' '26.883.59'
Set-CellText $ws 'E20' '  -0.77%  '
Set-CellText $ws 'D21' '5.011'
Set-CellText $ws 'E21' '  -0.44%  '
Set-CellText $ws 'D22' '2.119.29'
Set-CellText $ws 'E22' '  +1.63%  '
Set-CellText $ws 'D23' '10.36'
Set-CellText $ws 'E23' '  -1.40%  '
Set-CellText $ws 'D24' '6.372'
Set-CellText $ws 'E24' '  -1.34%  '
Set-CellText $ws 'D25' '148.09'
Set-CellText $ws 'D26' '1.783'
Set-CellText $ws 'E26' '  -3.15%  '
Set-CellText $ws 'D27' '17.86'
Set-CellText $ws 'E27' '  -0.62%  '
Set-CellText $ws 'D28' '2.089'
Set-CellText $ws 'E28' '  -1.02%  '
Set-CellText $ws 'D29' '113.38'
Set-CellText $ws 'E29' '  +0.62%  '
Set-CellText $ws 'D30' '4.704'
Set-CellText $ws 'E30' '  +1.00%  '
Set-CellText $ws 'D31' '4.713'
Set-CellText $ws 'E31' '  +0.22%  '
Set-CellText $ws 'D32' '0.09131'
Set-CellText $ws 'E32' '  +0.19%  '
Set-CellText $ws 'D33' '0.05072'
Set-CellText $ws 'E33' '  -1.07%  '
Set-CellText $ws 'B34' 'ImmutableX'
Set-CellText $ws 'C34' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-CellText $ws 'D34' '0.7469'
Set-CellText $ws 'E34' '  +2.93%  '
Set-CellText $ws 'B35' 'HuobiToken'
Set-CellText $ws 'C35' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-CellText $ws 'D35' '2.982'
Set-CellText $ws 'E35' '  -2.79%  '
Set-CellText $ws 'D36' '1.155'
Set-CellText $ws 'E36' '  -0.24%  '
Set-CellText $ws 'D37' '3.227'
Set-CellText $ws 'E37' '  +6.25%  '
Set-CellText $ws 'D38' '2.525'
Set-CellText $ws 'E38' '  +0.82%  '
Set-CellText $ws 'D39' '0.5617'
Set-CellText $ws 'E39' '  +5.23%  '
Set-CellText $ws 'D40' '0.01997'
Set-CellText $ws 'E40' '  -2.13%  '
Set-CellText $ws 'D41' '1.075'
Set-CellText $ws 'E41' '  +0.11%  '
Set-CellText $ws 'D42' '6.617'
Set-CellText $ws 'E42' '  +0.78%  '
Set-CellText $ws 'D43' '115.28'
Set-CellText $ws 'E43' '  -0.47%  '
Set-CellText $ws 'D44' '8.564'
Set-CellText $ws 'E44' '  +3.20%  '
Set-CellText $ws 'D45' '0.1475'
Set-CellText $ws 'E45' '  +0.58%  '
Set-CellText $ws 'D46' '0.4732'
Set-CellText $ws 'E46' '  +1.97%  '
Set-CellText $ws 'D47' '1.002'
Set-CellText $ws 'E47' '  +0.18%  '
Set-CellText $ws 'D48' '10.12'
Set-CellText $ws 'E48' '  +1.40%  '
Set-CellText $ws 'D49' '1.566'
Set-CellText $ws 'E49' '  -0.42%  '
Set-CellText $ws 'D50' '36.88'
Set-CellText $ws 'E50' '  +0.91%  '
Set-CellText $ws 'D51' '63.05'
Set-CellText $ws 'E51' '  -1.23%  '
